$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Capture the "old" look (style index 15) that currently lives on row 493
#    (C and E:L) and stamp it onto row 503 BEFORE we repaint rows 493-496,
#    since row 503 ends up using that same look in the final sheet.
# ---------------------------------------------------------------------------
$ws.Range("C493").Copy()
$ws.Range("C503").PasteSpecial(-4122)
$ws.Range("E493:L493").Copy()
$ws.Range("E503:L503").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Re-paint rows 493-496 (C, E:L) with the look used by row 492 (style 1),
#    leaving column D untouched.
# ---------------------------------------------------------------------------
$ws.Range("C492").Copy()
$ws.Range("C493:C496").PasteSpecial(-4122)
$ws.Range("E492:L492").Copy()
$ws.Range("E493:L493").PasteSpecial(-4122)
$ws.Range("E493:L493").Copy()
$ws.Range("E494:L496").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Give the brand-new rows 497-502 the same look as row 492 (style 1).
# ---------------------------------------------------------------------------
$ws.Range("C492").Copy()
$ws.Range("C497:C502").PasteSpecial(-4122)
$ws.Range("E492:L492").Copy()
$ws.Range("E497:L497").PasteSpecial(-4122)
$ws.Range("E497:L497").Copy()
$ws.Range("E498:L502").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Fill in the newly reported data (14 Aug and the days around it).
# ---------------------------------------------------------------------------
$data = @{
    497 = @(28716.0, 563.0, 65.0, 64.0, 40.0, 155901.0, 60262.0, 5941.0, 222104.0)
    498 = @(28730.0, 564.0, 61.0, 59.0, 39.0, 156053.0, 60289.0, 5941.0, 222283.0)
    499 = @(28771.0, 565.0, 68.0, 65.0, 41.0, 157889.0, 60963.0, 5942.0, 224794.0)
    500 = @(28865.0, 565.0, 62.0, 61.0, 36.0, 158919.0, 61676.0, 5942.0, 226537.0)
    501 = @(28904.0, 566.0, 59.0, 58.0, 32.0, 160557.0, 62633.0, 5943.0, 229133.0)
    502 = @(28956.0, 566.0, 61.0, 60.0, 33.0, 162744.0, 63766.0, 5943.0, 232453.0)
    503 = @(28995.0, 566.0, 58.0, 56.0, 34.0, 165134.0, 64846.0, 6943.0, 235923.0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value  = $vals[0]   # C - confirmados
    $ws.Cells.Item($row, 5).Value  = $vals[1]   # E - mortes
    $ws.Cells.Item($row, 6).Value  = $vals[2]   # F - total internacoes
    $ws.Cells.Item($row, 7).Value  = $vals[3]   # G - internacoes confirmadas
    $ws.Cells.Item($row, 8).Value  = $vals[4]   # H - UTI (total)
    $ws.Cells.Item($row, 9).Value  = $vals[5]   # I - 1a dose
    $ws.Cells.Item($row, 10).Value = $vals[6]   # J - 2a dose
    $ws.Cells.Item($row, 11).Value = $vals[7]   # K - dose unica
    $ws.Cells.Item($row, 12).Value = $vals[8]   # L - total vac
}

# ---------------------------------------------------------------------------
# 5. Extend the "novos" (D) shared formula down through the new rows.
# ---------------------------------------------------------------------------
$ws.Range("D497:D503").Formula = "=(C497-C496)"

# ---------------------------------------------------------------------------
# 6. The workbook/sheet filter view got a new internal id on re-export;
#    rename the matching defined name to track it (the one piece of that
#    identifier the Excel object model lets us touch).
# ---------------------------------------------------------------------------
$wb.Names.Item(1).Name = "Z_12621F90_3527_4EC5_8665_E3613A458E47_.wvu.FilterData"
